# "CSV fil upload was integrated"
# - Replace the placeholder REG NO value in A2 with the real uploaded value.
# - Remove the stray leftover row 3 (an orphaned LEADERSHIP value with no
#   matching student record), shrinking the used range back to A1:G2.
# - Leave the active selection on H2, just past the data, as it was after
#   the edit was made in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student registration number that was uploaded via CSV.
$ws.Range("A2").Value = "2021C123450003"

# Drop the orphaned row 3 (only G3 had a value) left over from the old data.
$ws.Range("A3:G3").ClearContents()

# Match the final cursor position recorded in the saved workbook.
$ws.Range("H2").Select()
